$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.700.40'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.599.53'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.07'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.0618'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.61'
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.823.51'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.603.00'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.92'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.681.96'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '208.65'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.80'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('E23').Value = '  -3.34%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.89'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.05'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.23'
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('E28').Value = '  +2.06%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.28'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  +0.96%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.659'
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.94'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.293.51'
$ws.Range('E35').Value = '  -1.60%  '
$ws.Range('E36').Value = '  -2.32%  '
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('E39').Value = '  +2.31%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.41'
$ws.Range('E41').Value = '  +1.41%  '
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.787'
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.62'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.735.95'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.900'
$ws.Range('E46').Value = '  +8.50%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.08'
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('E49').Value = '  +1.80%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0506'
$ws.Range('E50').Value = '  -0.26%  '
$ws.Range('E51').Value = '  +0.55%  '
